# Update "想去人数" (want-to-go count) figures that changed between scrapes.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 5291
$ws1.Range("F5").Value = 5291
$ws1.Range("F9").Value = 8849
$ws1.Range("F28").Value = 7159
$ws1.Range("F32").Value = 41
$ws1.Range("F41").Value = 2551
$ws1.Range("F50").Value = 1135

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 21

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 5291
$ws4.Range("F4").Value = 5291
$ws4.Range("F8").Value = 8849
$ws4.Range("F9").Value = 8849
$ws4.Range("F13").Value = 21
$ws4.Range("F31").Value = 7159
$ws4.Range("F40").Value = 2551
$ws4.Range("F50").Value = 1135
